# Update countries & provincias Spain
# Applies the 6-Aug-2020 17:27 -> 18:44 data refresh to the "Pais" sheet:
#  - Swap the Zambia/Luxemburgo pair (rows 95-96) and the
#    Timor Oriental/Santa Lucia pair (rows 202-203) to reflect the new sort
#    order produced by the refreshed case counts.
#  - Update the updated-at timestamp banner in A1.
#  - Update the numeric case/death/recovery counters for the affected
#    countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Country name swaps (re-sort ties) -------------------------------------
# Row 95 used to be Zambia, row 96 used to be Luxemburgo; after the refresh
# Luxemburgo overtakes Zambia in the ranking, so the labels swap while the
# (now-updated) numeric columns follow below.
$ws.Range("A95").Value = "Luxemburgo"
$ws.Range("A96").Value = "Zambia"

# Row 202 used to be Timor Oriental, row 203 used to be Santa Lucia; these two
# are tied on every numeric column, so only the labels need to swap.
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- Numeric data refresh ----------------------------------------------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 4991802
$ws.Range("C4").Value = 18234
$ws.Range("D4").Value = 2542849
$ws.Range("E4").Value = 2287017
$ws.Range("G4").Value = 335
$ws.Range("H4").Value = 161936

# Row 5: Brasil
$ws.Range("B5").Value = 2873304
$ws.Range("C5").Value = 10543
$ws.Range("E5").Value = 754975
$ws.Range("G5").Value = 274
$ws.Range("H5").Value = 97692

# Row 6: India
$ws.Range("B6").Value = 2021407
$ws.Range("C6").Value = 58168
$ws.Range("D6").Value = 1374420
$ws.Range("E6").Value = 605360
$ws.Range("G6").Value = 888
$ws.Range("H6").Value = 41627

# Row 11: Chile
$ws.Range("B11").Value = 366671
$ws.Range("C11").Value = 1948
$ws.Range("D11").Value = 340168
$ws.Range("E11").Value = 16614
$ws.Range("G11").Value = 97
$ws.Range("H11").Value = 9889

# Row 15: Reino Unido
$ws.Range("G15").Value = 49
$ws.Range("H15").Value = 46413

# Row 19: Italia
$ws.Range("B19").Value = 249204
$ws.Range("C19").Value = 401
$ws.Range("D19").Value = 201323
$ws.Range("E19").Value = 12694
$ws.Range("G19").Value = 6
$ws.Range("H19").Value = 35187

# Row 38: Republica Dominicana
$ws.Range("B38").Value = 76536
$ws.Range("C38").Value = 876
$ws.Range("D38").Value = 40539
$ws.Range("E38").Value = 34751
$ws.Range("G38").Value = 24
$ws.Range("H38").Value = 1246

# Row 43: Emiratos Arabes Unidos
$ws.Range("B43").Value = 61845
$ws.Range("C43").Value = 239
$ws.Range("D43").Value = 55739
$ws.Range("E43").Value = 5752
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 354

# Row 47: Guatemala
$ws.Range("B47").Value = 54339
$ws.Range("C47").Value = 830
$ws.Range("D47").Value = 42070
$ws.Range("E47").Value = 10150
$ws.Range("G47").Value = 47
$ws.Range("H47").Value = 2119

# Row 75: Chequia
$ws.Range("B75").Value = 17648
$ws.Range("C75").Value = 119
$ws.Range("D75").Value = 12196
$ws.Range("E75").Value = 5063
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 389

# Row 80: Bosnia y Herzegovina
$ws.Range("B80").Value = 13396
$ws.Range("C80").Value = 258
$ws.Range("D80").Value = 7042
$ws.Range("E80").Value = 5970

# Row 95: Zambia (now showing Luxemburgo data)
$ws.Range("B95").Value = 7073
$ws.Range("C95").Value = 66
$ws.Range("D95").Value = 5750
$ws.Range("E95").Value = 1204
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 119

# Row 96: Luxemburgo (now showing Zambia data)
$ws.Range("B96").Value = 7022
$ws.Range("D96").Value = 5667
$ws.Range("E96").Value = 1179
$ws.Range("H96").Value = 176

# Row 100: Libano
$ws.Range("B100").Value = 5672
$ws.Range("C100").Value = 255
$ws.Range("D100").Value = 1974
$ws.Range("E100").Value = 3628
$ws.Range("G100").Value = 2
$ws.Range("H100").Value = 70

# Row 103: Grecia
$ws.Range("B103").Value = 5123
$ws.Range("C103").Value = 149
$ws.Range("E103").Value = 3539

# Row 112: Hong Kong
$ws.Range("E112").Value = 1346
$ws.Range("G112").Value = 3
$ws.Range("H112").Value = 46

# Row 142: Jordania
$ws.Range("B142").Value = 1232
$ws.Range("C142").Value = 1
$ws.Range("D142").Value = 1171
$ws.Range("E142").Value = 50

# Row 163: Reunion
$ws.Range("B163").Value = 671
$ws.Range("C163").Value = 1
$ws.Range("E163").Value = 74

# Row 175: Martinica
$ws.Range("B175").Value = 276
$ws.Range("C175").Value = 7
$ws.Range("E175").Value = 163

# --- Updated-at banner -------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Agosto de 2020 a las 18:44"
